# Generated edit script for Weight_matrix.xlsx change
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: new literal weight values (C11:O11) ---
$ws.Range("C11").Value = 0.3839762433259343
$ws.Range("D11").Value = 0.805702416918429
$ws.Range("E11").Value = 0.7237928304873912
$ws.Range("F11").Value = 0.8615560640732265
$ws.Range("G11").Value = 0.7213456553589541
$ws.Range("H11").Value = 0.7922391385072409
$ws.Range("I11").Value = 0.9354720841859105
$ws.Range("J11").Value = 0.8819760231500618
$ws.Range("K11").Value = 0.9341068301225918
$ws.Range("L11").Value = 0.8950496434065165
$ws.Range("M11").Value = 0.9259982638888886
$ws.Range("N11").Value = 0.9715391621129326
$ws.Range("O11").Value = 1.0

# --- Row 12: new literal weight values (C12:O12) ---
$ws.Range("C12").Value = 0.36910994764397914
$ws.Range("D12").Value = 0.6409090909090911
$ws.Range("E12").Value = 0.2183789364997419
$ws.Range("F12").Value = 0.6175182481751827
$ws.Range("G12").Value = 0.2183789364997419
$ws.Range("H12").Value = 0.4976470588235295
$ws.Range("I12").Value = 0.5395408163265307
$ws.Range("J12").Value = 0.714527027027027
$ws.Range("K12").Value = 0.6477794793261868
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.94
$ws.Range("N12").Value = 0.5818431911966988
$ws.Range("O12").Value = 0.9276315789473686

# --- Row 13: formula changes from average (C12+C11)/2 to product C12*C11 ---
$ws.Range("C13").Formula = "=C12*C11"
$ws.Range("D13").Formula = "=D12*D11"
$ws.Range("E13").Formula = "=E12*E11"
$ws.Range("F13").Formula = "=F12*F11"
$ws.Range("G13").Formula = "=G12*G11"
$ws.Range("H13").Formula = "=H12*H11"
$ws.Range("I13").Formula = "=I12*I11"
$ws.Range("J13").Formula = "=J12*J11"
$ws.Range("K13").Formula = "=K12*K11"
$ws.Range("L13").Formula = "=L12*L11"
$ws.Range("M13").Formula = "=M12*M11"
$ws.Range("N13").Formula = "=N12*N11"
$ws.Range("O13").Formula = "=O12*O11"

# --- Row 26: new literal values (C26:O26) ---
$ws.Range("C26").Value = 0.9048991354466858
$ws.Range("D26").Value = 0.4617647058823529
$ws.Range("E26").Value = 1.0
$ws.Range("F26").Value = 0.48307692307692307
$ws.Range("G26").Value = 0.9486404833836857
$ws.Range("H26").Value = 0.02390377588306943
$ws.Range("I26").Value = 0.00679242017824695
$ws.Range("J26").Value = 0.6061776061776062
$ws.Range("K26").Value = 0.3588571428571428
$ws.Range("L26").Value = 0.18416422287390027
$ws.Range("M26").Value = 0.7302325581395348
$ws.Range("N26").Value = 0.591337099811676
$ws.Range("O26").Value = 0.027689594356261022

# --- Row 28: formulas replace the hardcoded 1's, now referencing row 13 ---
$ws.Range("C28").Formula = "=C13"
$ws.Range("D28").Formula = "=D13"
$ws.Range("E28").Formula = "=E13"
$ws.Range("F28").Formula = "=F13"
$ws.Range("G28").Formula = "=G13"
$ws.Range("H28").Formula = "=H13"
$ws.Range("I28").Formula = "=I13"
$ws.Range("J28").Formula = "=J13"
$ws.Range("K28").Formula = "=K13"
$ws.Range("L28").Formula = "=L13"
$ws.Range("M28").Formula = "=M13"
$ws.Range("N28").Formula = "=N13"
$ws.Range("O28").Formula = "=O13"

# --- Update the active selection to M21 (was M19) ---
$ws.Range("M21").Select()

